# Replace the "thermochemical water splitting" hydrogen production pathway
# with "hydrocarbon partial oxidation" on the RHPF sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("RHPF")

$old = "thermochemical water splitting"
$new = "hydrocarbon partial oxidation"

# Column header (row 1) and row header (column A) both reference the pathway name.
$ws.Range("F1").Value = $new
$ws.Range("A6").Value = $new
